# Zapojeni.xlsx - "Code ready for hand in"
# Adds APB2/APB1 bus info block (C8:E9) and a timer/adc clock-divider
# calculation block (F16:I19) to List1, then updates the sheet view
# (zoom + active selection) to match where the author left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- APB2 / APB1 info block -------------------------------------------------
$ws.Range("C8").Value = "APB2"
$ws.Range("D8").Value = 96
$ws.Range("E8").Value = "MHZ"
$ws.Range("C9").Value = "APB1"

# --- timer / adc clock calculation block ------------------------------------
$ws.Range("F16").Value = "timer"
$ws.Range("G16").Value = 500000
$ws.Range("F17").Value = "adc"
$ws.Range("G17").Value = 4096

$ws.Range("I16").Formula = "=G17*122"
$ws.Range("F19").Formula = "=G16/G17"

# --- sheet view: zoom + final selection -------------------------------------
$win = $excel.ActiveWindow
$win.Zoom = 70
$ws.Range("G17").Select()

$wb.Save()
